$wb = $excel.ActiveWorkbook

# The three new country tabs were cloned from the "Italy" sheet (same
# column widths / cell styling) and dropped in after "Greece", which
# stays the last of the original eleven tabs.
$italy  = $wb.Worksheets.Item("Italy")
$greece = $wb.Worksheets.Item("Greece")

# --- Netherlands -----------------------------------------------------
$italy.Copy($null, $greece)
$netherlands = $wb.Worksheets.Item($wb.Worksheets.Count)
$netherlands.Name = "Netherlands"
$netherlands.Range("B4").Value = "NGC-3145/T2160"
$netherlands.Range("B2").Value = "Netherlands Market"
$netherlands.Range("B4").Select()

# --- Austria -----------------------------------------------------------
$netherlands.Copy($null, $netherlands)
$austria = $wb.Worksheets.Item($wb.Worksheets.Count)
$austria.Name = "Austria"
$austria.Range("B4").Value = "NGC-4320/T2276"
$austria.Range("B2").Value = "Austria Market"
$austria.Range("B4").Select()

# --- Denmark -----------------------------------------------------------
$austria.Copy($null, $austria)
$denmark = $wb.Worksheets.Item($wb.Worksheets.Count)
$denmark.Name = "Denmark"
$denmark.Range("B4").Value = "NGC-2913/T2756"
$denmark.Range("B2").Value = "Denmark Market"
$denmark.Range("B4").Select()

# Netherlands is left as the active/selected tab once all three sheets
# have been added.
$netherlands.Select()
